$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.343.47"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "'1.592.71"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'211.64"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'19.37"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'1.816.57"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.619.70"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "'64.55"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'26.356.04"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").Value = "'212.55"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "'144.76"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").Value = "'15.20"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'0.0501"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("D34").Value = "'1.337.95"
$ws.Range("E34").Value = "  +4.14%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").Value = "'0.601"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "'1.48"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  -18.51%  "
$ws.Range("D41").Value = "'5.81"
$ws.Range("E41").Value = "  +5.16%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "'0.763"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "'1.729.19"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'61.60"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "'87.76"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").Value = "'0.0983"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  -0.63%  "
